$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F column (Baseline Start) for group header rows
$ws.Range("F3").Value = 42116.6117988432
$ws.Range("F4").Value = 42116.6117988432
$ws.Range("F11").Value = 42116.6117988432
$ws.Range("F17").Value = 42116.6117988432
$ws.Range("F22").Value = 42116.6117988432
$ws.Range("F26").Value = 42116.6117988432
$ws.Range("F33").Value = 42116.6117988432
$ws.Range("F35").Value = 42116.6117988432
$ws.Range("F43").Value = 42116.6117988432
$ws.Range("F51").Value = 42116.6117988432

# Update G column (Baseline End) for all rows 3-61
$ws.Range("G3").Value = 42126.6117988434
$ws.Range("G4").Value = 42126.6117988434
$ws.Range("G5").Value = 42126.6117988434
$ws.Range("G6").Value = 42126.6117988434
$ws.Range("G7").Value = 42126.6117988434
$ws.Range("G8").Value = 42126.6117988434
$ws.Range("G9").Value = 42126.6117988434
$ws.Range("G10").Value = 42126.6117988434
$ws.Range("G11").Value = 42126.6117988434
$ws.Range("G12").Value = 42126.6117988434
$ws.Range("G13").Value = 42126.6117988434
$ws.Range("G14").Value = 42126.6117988434
$ws.Range("G15").Value = 42126.6117988434
$ws.Range("G16").Value = 42126.6117988434
$ws.Range("G17").Value = 42126.6117988434
$ws.Range("G18").Value = 42126.6117988434
$ws.Range("G19").Value = 42126.6117988434
$ws.Range("G20").Value = 42126.6117988434
$ws.Range("G21").Value = 42126.6117988434
$ws.Range("G22").Value = 42126.6117988434
$ws.Range("G23").Value = 42126.6117988434
$ws.Range("G24").Value = 42126.6117988434
$ws.Range("G25").Value = 42126.6117988434
$ws.Range("G26").Value = 42126.6117988434
$ws.Range("G27").Value = 42126.6117988434
$ws.Range("G28").Value = 42126.6117988434
$ws.Range("G29").Value = 42126.6117988434
$ws.Range("G30").Value = 42126.6117988434
$ws.Range("G31").Value = 42126.6117988434
$ws.Range("G32").Value = 42126.6117988434
$ws.Range("G33").Value = 42126.6117988434
$ws.Range("G34").Value = 42126.6117988434
$ws.Range("G35").Value = 42126.6117988434
$ws.Range("G36").Value = 42126.6117988434
$ws.Range("G37").Value = 42126.6117988434
$ws.Range("G38").Value = 42126.6117988434
$ws.Range("G39").Value = 42126.6117988434
$ws.Range("G40").Value = 42126.6117988434
$ws.Range("G41").Value = 42126.6117988434
$ws.Range("G42").Value = 42126.6117988434
$ws.Range("G43").Value = 42126.6117988434
$ws.Range("G44").Value = 42126.6117988434
$ws.Range("G45").Value = 42126.6117988434
$ws.Range("G46").Value = 42126.6117988434
$ws.Range("G47").Value = 42126.6117988434
$ws.Range("G48").Value = 42126.6117988434
$ws.Range("G49").Value = 42126.6117988434
$ws.Range("G50").Value = 42126.6117988434
$ws.Range("G51").Value = 42126.6117988434
$ws.Range("G52").Value = 42126.6117988434
$ws.Range("G53").Value = 42126.6117988434
$ws.Range("G54").Value = 42126.6117988434
$ws.Range("G55").Value = 42126.6117988434
$ws.Range("G56").Value = 42126.6117988434
$ws.Range("G57").Value = 42126.6117988434
$ws.Range("G58").Value = 42126.6117988434
$ws.Range("G59").Value = 42126.6117988434
$ws.Range("G60").Value = 42126.6117988434
$ws.Range("G61").Value = 42126.6117988434

# Update K column (Fixed Cost) for group header rows
$ws.Range("K3").Value = 2700210.30574707
$ws.Range("K4").Value = 2700210.30574707
$ws.Range("K11").Value = 2700210.30574707
$ws.Range("K17").Value = 2700210.30574707
$ws.Range("K22").Value = 2700210.30574707
$ws.Range("K26").Value = 2700210.30574707
$ws.Range("K33").Value = 2700210.30574707
$ws.Range("K35").Value = 2700210.30574707
$ws.Range("K43").Value = 2700210.30574707
$ws.Range("K51").Value = 2700210.30574707

# Update N column (Total Cost) for group header rows
$ws.Range("N3").Value = 2952108.70818848
$ws.Range("N4").Value = 2952108.70818848
$ws.Range("N11").Value = 2952108.70818848
$ws.Range("N17").Value = 2952108.70818848
$ws.Range("N22").Value = 2952108.70818848
$ws.Range("N26").Value = 2952108.70818848
$ws.Range("N33").Value = 2952108.70818848
$ws.Range("N35").Value = 2952108.70818848
$ws.Range("N43").Value = 2952108.70818848
$ws.Range("N51").Value = 2952108.70818848

